$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column A; this shifts the existing U-233..Th-229
# data (and their column widths) one column to the right, preserving their
# stored widths exactly.
$ws.Columns.Item(1).Insert()

# New header for the inserted "Lab. #" column.
$ws.Range("A1").Value = "Lab. #"

# "Lab. #" values for the 11 data rows (rows 2-12).
$labValues = @(10815, 11069, 10815, 11070, 10815, 11071, 10815, 11072, 10815, 11074, 10815)
for ($i = 0; $i -lt $labValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $labValues[$i]
}

# New column's width.
$ws.Columns.Item(1).ColumnWidth = 7.7109375

# Shade every other data row (2,4,6,8,10,12) with a light-green fill, across
# the full A:I extent of that row's data.
$shadeRows = @(2, 4, 6, 8, 10, 12)
foreach ($r in $shadeRows) {
    $rng = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 9))
    $rng.Interior.Color = 12379352
}
